# Hjemme passive tweaks lichtwark deleted values
# Updates the row-1 header values and the row-2/row-3 measurement values
# for columns B:E, then restores the visible selection to B1:E3 (matching
# the saved "before close" selection captured in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header/id) values -------------------------------------------
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 values ---------------------------------------------------------
$ws.Range("B2").Value = 2.9836900550838363
$ws.Range("C2").Value = 6.317758840049482
$ws.Range("D2").Value = 3.1804347241046571
$ws.Range("E2").Value = 3.6586700155965675

# --- Row 3 values ---------------------------------------------------------
$ws.Range("B3").Value = 7.900352088866569
$ws.Range("C3").Value = 10.000243686011228
$ws.Range("D3").Value = 5.0392132211217087
$ws.Range("E3").Value = 1.7873213578991689

# --- Selection: narrow the saved selection from B1:AY3 to B1:E3 ----------
$ws.Activate() | Out-Null
$ws.Range("B1:E3").Select() | Out-Null
